# Updates the cryptos list on the active worksheet to reflect the latest
# scraped prices / 1h volume changes, including re-ranking a few coins
# whose relative order changed (rows 12/13 and 41/42/43).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value into a cell while guaranteeing it stays plain text
# (many "Price" values look like numbers, e.g. "98.31", and Excel would
# otherwise silently convert them to the Number type on assignment).
function Set-TextValue($addr, $value) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

# --- Row 2: Bitcoin ---------------------------------------------------
Set-TextValue "D2" "42.975.29"
Set-TextValue "E2" "  +0.13%  "

# --- Row 3: Ethereum ---------------------------------------------------
Set-TextValue "D3" "2.546.92"
Set-TextValue "E3" "  +0.11%  "

# --- Row 4: TetherUSD ---------------------------------------------------
Set-TextValue "E4" "  -0.05%  "

# --- Row 5: BNB ---------------------------------------------------
Set-TextValue "D5" "304.49"
Set-TextValue "E5" "  +1.91%  "

# --- Row 6: Solana ---------------------------------------------------
Set-TextValue "D6" "98.31"
Set-TextValue "E6" "  +4.43%  "

# --- Row 7: XRP ---------------------------------------------------
Set-TextValue "E7" "  +0.92%  "

# --- Row 8: USDC ---------------------------------------------------
Set-TextValue "E8" "  -0.05%  "

# --- Row 9: Cardano ---------------------------------------------------
Set-TextValue "D9" "0.546"
Set-TextValue "E9" "  -0.40%  "

# --- Row 10: Avalanche ---------------------------------------------------
Set-TextValue "D10" "37.11"
Set-TextValue "E10" "  +3.94%  "

# --- Row 11: Dogecoin ---------------------------------------------------
Set-TextValue "E11" "  +3.01%  "

# --- Row 12 / 13: TRON and Polkadot swap rank ---------------------------
Set-TextValue "B12" "Polkadot"
Set-TextValue "C12" "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
Set-TextValue "D12" "7.78"
Set-TextValue "E12" "  +1.80%  "

Set-TextValue "B13" "TRON"
Set-TextValue "C13" "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue "D13" "0.116"
Set-TextValue "E13" "  +3.25%  "

# --- Row 14: WrappedliquidstakedEther2.0 ---------------------------------
Set-TextValue "D14" "2.939.88"
Set-TextValue "E14" "  +0.23%  "

# --- Row 15: WrappedEther ---------------------------------------------------
Set-TextValue "D15" "2.554.28"
Set-TextValue "E15" "  +0.38%  "

# --- Row 16: Chainlink ---------------------------------------------------
Set-TextValue "D16" "15.10"
Set-TextValue "E16" "  +7.42%  "

# --- Row 17: Polygon ---------------------------------------------------
Set-TextValue "D17" "0.871"
Set-TextValue "E17" "  +0.33%  "

# --- Row 18: WrappedBTC ---------------------------------------------------
Set-TextValue "D18" "42.968.89"
Set-TextValue "E18" "  +0.03%  "

# --- Row 19: InternetComputer(DFINITY) -----------------------------------
Set-TextValue "D19" "13.88"
Set-TextValue "E19" "  +7.21%  "

# --- Row 20: ShibaInu ---------------------------------------------------
Set-TextValue "E20" "  +1.30%  "

# --- Row 21: Uniswap ---------------------------------------------------
Set-TextValue "D21" "6.58"
Set-TextValue "E21" "  +0.07%  "

# --- Row 22: Litecoin ---------------------------------------------------
Set-TextValue "D22" "71.98"
Set-TextValue "E22" "  +0.51%  "

# --- Row 23: BitcoinCash ---------------------------------------------------
Set-TextValue "D23" "254.95"
Set-TextValue "E23" "  -0.30%  "

# --- Row 24: PancakeSwap ---------------------------------------------------
Set-TextValue "E24" "  +2.30%  "

# --- Row 25: ImmutableX ---------------------------------------------------
Set-TextValue "E25" "  -1.43%  "

# --- Row 26: EthereumClassic ---------------------------------------------
Set-TextValue "D26" "28.07"
Set-TextValue "E26" "  -3.65%  "

# --- Row 27: Dai ---------------------------------------------------
Set-TextValue "D27" "0.999"
Set-TextValue "E27" "  -0.08%  "

# --- Row 28: Cosmos ---------------------------------------------------
Set-TextValue "D28" "10.27"
Set-TextValue "E28" "  +2.84%  "

# --- Row 29: InjectiveProtocol ---------------------------------------------
Set-TextValue "D29" "37.78"
Set-TextValue "E29" "  +1.82%  "

# --- Row 30: Toncoin ---------------------------------------------------
Set-TextValue "E30" "  -1.77%  "

# --- Row 31: Filecoin ---------------------------------------------------
Set-TextValue "D31" "6.17"
Set-TextValue "E31" "  +4.44%  "

# --- Row 32: Monero ---------------------------------------------------
Set-TextValue "D32" "158.59"
Set-TextValue "E32" "  +3.85%  "

# --- Row 33: Celestia ---------------------------------------------------
Set-TextValue "D33" "19.47"
Set-TextValue "E33" "  +14.94%  "

# --- Row 34: ARBITRUM ---------------------------------------------------
Set-TextValue "D34" "2.15"
Set-TextValue "E34" "  -0.39%  "

# --- Row 35: Hedera ---------------------------------------------------
Set-TextValue "D35" "0.0803"
Set-TextValue "E35" "  +1.25%  "

# --- Row 36: LidoDAOToken ---------------------------------------------------
Set-TextValue "D36" "3.31"
Set-TextValue "E36" "  -2.03%  "

# --- Row 37: WEMIXToken ---------------------------------------------------
Set-TextValue "E37" "  -4.32%  "

# --- Row 38: Kaspa ---------------------------------------------------
Set-TextValue "E38" "  +2.11%  "

# --- Row 39: EnergySwap ---------------------------------------------------
Set-TextValue "D39" "25.40"
Set-TextValue "E39" "  +9.91%  "

# --- Row 40: Stellar ---------------------------------------------------
Set-TextValue "E40" "  -0.05%  "

# --- Row 41 / 42 / 43: ApeXProtocol, NEARProtocol, RenderToken re-rank ---
Set-TextValue "B41" "NEARProtocol"
Set-TextValue "C41" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue "D41" "3.44"
Set-TextValue "E41" "  +0.71%  "

Set-TextValue "B42" "RenderToken"
Set-TextValue "C42" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D42" "3.91"
Set-TextValue "E42" "  +0.49%  "

Set-TextValue "B43" "ApeXProtocol"
Set-TextValue "C43" "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
Set-TextValue "D43" "2.08"
Set-TextValue "E43" "  +30.49%  "

# --- Row 44: Maker ---------------------------------------------------
Set-TextValue "D44" "2.099.50"
Set-TextValue "E44" "  +0.66%  "

# --- Row 45: VeChain ---------------------------------------------------
Set-TextValue "D45" "0.0306"
Set-TextValue "E45" "  -1.18%  "

# --- Row 46: FirstDigitalUSD ---------------------------------------------
Set-TextValue "E46" "  -0.11%  "

# --- Row 47: BitcoinSV ---------------------------------------------------
Set-TextValue "D47" "86.73"
Set-TextValue "E47" "  +3.06%  "

# --- Row 48: FraxShare ---------------------------------------------------
Set-TextValue "D48" "8.91"
Set-TextValue "E48" "  +0.23%  "

# --- Row 49: ordi ---------------------------------------------------
Set-TextValue "D49" "75.46"
Set-TextValue "E49" "  +10.05%  "

# --- Row 50: RocketPoolETH ---------------------------------------------------
Set-TextValue "D50" "2.798.15"
Set-TextValue "E50" "  +0.18%  "

# --- Row 51: Algorand ---------------------------------------------------
Set-TextValue "E51" "  +3.18%  "
